# Weekly fruit/vegetable price update: insert two new price records
# (rows 424-425) for "Zapallo italiano" / Femacal de La Calera, pushing
# the existing rows 424-431 down to 426-433.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 424 - this shifts the
# old rows 424..431 down to 426..433, preserving their data untouched.
$ws.Rows("424:425").Insert()

# --- Row 424: new record -------------------------------------------------
$ws.Cells.Item(424, 1).Value = 3
$ws.Cells.Item(424, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(424, 3).Value = "Coquimbo"
$ws.Cells.Item(424, 4).Value = 44656
$ws.Cells.Item(424, 5).Value = 5
$ws.Cells.Item(424, 6).Value = 100112032
$ws.Cells.Item(424, 7).Value = "Zapallo italiano"
$ws.Cells.Item(424, 8).Value = "Sin especificar"
$ws.Cells.Item(424, 9).Value = "Primera"
$ws.Cells.Item(424, 10).Value = 130
$ws.Cells.Item(424, 11).Value = 4500
$ws.Cells.Item(424, 12).Value = 5000
$ws.Cells.Item(424, 13).Value = 4731
$ws.Cells.Item(424, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(424, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(424, 16).Value = 131
$ws.Cells.Item(424, 17).Value = 36
$ws.Cells.Item(424, 18).Value = "Hortaliza"

# --- Row 425: new record -------------------------------------------------
$ws.Cells.Item(425, 1).Value = 3
$ws.Cells.Item(425, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(425, 3).Value = "Coquimbo"
$ws.Cells.Item(425, 4).Value = 44656
$ws.Cells.Item(425, 5).Value = 5
$ws.Cells.Item(425, 6).Value = 100112032
$ws.Cells.Item(425, 7).Value = "Zapallo italiano"
$ws.Cells.Item(425, 8).Value = "Sin especificar"
$ws.Cells.Item(425, 9).Value = "Primera"
$ws.Cells.Item(425, 10).Value = 175
$ws.Cells.Item(425, 11).Value = 9000
$ws.Cells.Item(425, 12).Value = 10000
$ws.Cells.Item(425, 13).Value = 9486
$ws.Cells.Item(425, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(425, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(425, 16).Value = 136
$ws.Cells.Item(425, 17).Value = 70
$ws.Cells.Item(425, 18).Value = "Hortaliza"
